# Commit: "refactored code to remove duplicates"
#
# The "customer" sheet's D2 cell used to hold a single reference to one
# address record ("reference:address@id#1"). It is refactored to a
# "listReference" that can point at multiple address records
# ("listReference:address@id#1,address@id#2") instead of duplicating a
# separate single-reference column.

$wb = $excel.ActiveWorkbook

$customer = $wb.Worksheets.Item("customer")
$customer.Range("D2").Value = "listReference:address@id#1,address@id#2"

# Column B ("name" / "customer1") picks up a best-fit width once the sheet
# is touched/resaved, matching the rest of the workbook's autosized columns.
$customer.Columns("B").AutoFit() | Out-Null

# The "address" sheet's active selection moved on to the next cell (F2)
# after the edit.
$address = $wb.Worksheets.Item("address")
$address.Range("F2").Select() | Out-Null
